# Weekly CompStat data refresh for the week of 6/10/2024 - 6/16/2024.
# Updates the report header (volume/week-of text) and all precinct crime-stat
# figures (counts + %-change) in rows 14-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (shared-string rich text runs): issue number + report week
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# ---------------------------------------------------------------------
# Cells that flip between the "no data" text placeholder ("0" / "***.*")
# and a real numeric entry (or vice-versa) need their style swapped too.
# Copy an existing cell that already has the right style (and, for the
# placeholder cells, the right text) onto the target, then set the value.
# ---------------------------------------------------------------------

# --> numeric style (s="15")
$ws.Range("C16").Copy($ws.Range("C15"))
$ws.Range("C16").Copy($ws.Range("D15"))
$ws.Range("C16").Copy($ws.Range("F15"))
$ws.Range("C16").Copy($ws.Range("G15"))
$ws.Range("C16").Copy($ws.Range("C22"))
$ws.Range("C16").Copy($ws.Range("C23"))
$ws.Range("C16").Copy($ws.Range("D23"))
$ws.Range("C16").Copy($ws.Range("F23"))
$ws.Range("C16").Copy($ws.Range("C27"))
$ws.Range("C16").Copy($ws.Range("D31"))
$ws.Range("C16").Copy($ws.Range("G31"))

# --> numeric style (s="16")
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("K14").Copy($ws.Range("H15"))
$ws.Range("K14").Copy($ws.Range("E23"))
$ws.Range("K14").Copy($ws.Range("E31"))
$ws.Range("K14").Copy($ws.Range("H31"))

# --> text placeholder "0" (style s="14")
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("C14").Copy($ws.Range("F30"))

# --> text placeholder "***.*" (style s="14")
$ws.Range("E14").Copy($ws.Range("E22"))

# Values for the cells that became numeric (the placeholder-text cells
# above already carry the correct text from the donor copy)
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("C22").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("C27").Value = 2
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100

# ---------------------------------------------------------------------
# Remaining cells: value-only updates (style unchanged)
# ---------------------------------------------------------------------
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = -93.103448275862
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = -53.846153846153
$ws.Range("N15").Value = -77.777777777777
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -36.363636363636
$ws.Range("I16").Value = 55
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = 7.843137254901
$ws.Range("L16").Value = 14.583333333333
$ws.Range("M16").Value = -38.888888888888
$ws.Range("N16").Value = -83.383685800604
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = -17.647058823529
$ws.Range("I17").Value = 89
$ws.Range("J17").Value = 99
$ws.Range("K17").Value = -10.101010101010
$ws.Range("L17").Value = -22.608695652173
$ws.Range("M17").Value = 23.611111111111
$ws.Range("N17").Value = -72.012578616352
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 133.333333333333
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = 37.931034482758
$ws.Range("L18").Value = -4.761904761904
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = -87.841945288753
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 33
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -5.714285714285
$ws.Range("I19").Value = 167
$ws.Range("J19").Value = 137
$ws.Range("K19").Value = 21.897810218978
$ws.Range("L19").Value = -4.571428571428
$ws.Range("M19").Value = 153.030303030303
$ws.Range("N19").Value = 8.441558441558
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 125
$ws.Range("J20").Value = 49
$ws.Range("K20").Value = -34.693877551020
$ws.Range("L20").Value = 3.225806451612
$ws.Range("N20").Value = -77.622377622377
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -27.272727272727
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 71
$ws.Range("H21").Value = 1.408450704225
$ws.Range("I21").Value = 391
$ws.Range("J21").Value = 368
$ws.Range("K21").Value = 6.25
$ws.Range("L21").Value = -8
$ws.Range("M21").Value = 30.333333333333
$ws.Range("N21").Value = -70.623591284748
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = -53.846153846153
$ws.Range("L22").Value = -53.846153846153
$ws.Range("M22").Value = 100
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = -55.555555555555
$ws.Range("L23").Value = -33.333333333333
$ws.Range("C24").Value = 17
$ws.Range("E24").Value = 6.25
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -19.444444444444
$ws.Range("I24").Value = 324
$ws.Range("J24").Value = 330
$ws.Range("K24").Value = -1.818181818181
$ws.Range("L24").Value = -40.983606557377
$ws.Range("M24").Value = 123.448275862069
$ws.Range("C25").Value = 12
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 24
$ws.Range("G25").Value = 26
$ws.Range("H25").Value = -7.692307692307
$ws.Range("I25").Value = 78
$ws.Range("J25").Value = 85
$ws.Range("K25").Value = -8.235294117647
$ws.Range("L25").Value = -68.924302788844
$ws.Range("C26").Value = 10
$ws.Range("E26").Value = 25
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 154
$ws.Range("J26").Value = 138
$ws.Range("K26").Value = 11.594202898550
$ws.Range("L26").Value = -13.483146067415
$ws.Range("M26").Value = -28.703703703703
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = 11.111111111111
$ws.Range("L27").Value = -37.5
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = -7.142857142857
$ws.Range("L28").Value = -45.833333333333
$ws.Range("L29").Value = -28.571428571428
$ws.Range("L30").Value = -42.857142857142
$ws.Range("J31").Value = 2
$ws.Range("K31").Value = -50
